$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.06602291937237484
$ws.Range("C2").Value = 0.001004469044575683
$ws.Range("B3").Value = 0.08041596580427651
$ws.Range("C3").Value = 0.001160914808582568
